# Daily TGP (Terminal Gate Pricing) refresh: roll each table's two
# displayed effective dates forward by one day (45954->45955, 45955->45958)
# and update the Diesel/ULP/PULP/e10 price columns (D:G) to the new day's
# published cents-per-litre figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A8" = 45958
    "D8" = 163.75
    "E8" = 159.53
    "F8" = 169.53
    "G8" = 159.69
    "A9" = 45958
    "D9" = 163.75
    "E9" = 159.53
    "F9" = 169.53
    "G9" = 159.69
    "A10" = 45958
    "D10" = 166.18
    "E10" = 161.81
    "F10" = 171.81
    "G10" = 162.28
    "A11" = 45955
    "D11" = 161.32
    "E11" = 158.15
    "F11" = 168.15
    "G11" = 158.31
    "A12" = 45955
    "D12" = 161.32
    "E12" = 158.15
    "F12" = 168.15
    "G12" = 158.31
    "A13" = 45955
    "D13" = 163.72
    "E13" = 160.49
    "F13" = 170.49
    "G13" = 160.96
    "A17" = 45958
    "D17" = 169.63
    "E17" = 164.59
    "F17" = 174.59
    "A18" = 45955
    "D18" = 167.16
    "E18" = 163.29
    "F18" = 173.29
    "A22" = 45958
    "D22" = 165.11
    "E22" = 160.81
    "F22" = 170.41
    "G22" = 161.99
    "A23" = 45958
    "D23" = 170.95
    "E23" = 165.53
    "F23" = 175.53
    "A24" = 45958
    "D24" = 170.76
    "E24" = 165.73
    "F24" = 175.73
    "A25" = 45958
    "D25" = 171.59
    "E25" = 165.12
    "F25" = 175.12
    "G25" = 164.95
    "A26" = 45958
    "D26" = 170.31
    "E26" = 166.68
    "F26" = 176.68
    "A27" = 45955
    "D27" = 162.44999999999999
    "E27" = 159.5
    "F27" = 169.1
    "G27" = 160.66999999999999
    "A28" = 45955
    "D28" = 168.49
    "E28" = 164.21
    "F28" = 174.21
    "A29" = 45955
    "D29" = 168.29
    "E29" = 164.41
    "F29" = 174.41
    "A30" = 45955
    "D30" = 169.12
    "E30" = 163.80000000000001
    "F30" = 173.8
    "G30" = 163.63
    "A31" = 45955
    "D31" = 167.85
    "E31" = 165.34
    "F31" = 175.34
    "A35" = 45958
    "D35" = 164.43
    "E35" = 159.02000000000001
    "F35" = 168.02
    "A36" = 45955
    "D36" = 161.97
    "E36" = 157.69999999999999
    "F36" = 166.71
    "A40" = 45958
    "D40" = 170.08
    "E40" = 164.3
    "F40" = 174.3
    "A41" = 45958
    "D41" = 169.79
    "E41" = 164.72
    "F41" = 174.72
    "A42" = 45955
    "D42" = 167.61
    "E42" = 163.01
    "F42" = 173.01
    "A43" = 45955
    "D43" = 167.33
    "E43" = 163.43
    "F43" = 173.43
    "A47" = 45958
    "D47" = 161.88
    "E47" = 159.86000000000001
    "F47" = 169.86
    "A48" = 45958
    "D48" = 161.87
    "E48" = 160.04
    "F48" = 170.04
    "A49" = 45955
    "D49" = 160.61000000000001
    "E49" = 159.08000000000001
    "F49" = 169.08
    "A50" = 45955
    "D50" = 160.59
    "E50" = 159.25
    "F50" = 169.25
    "A54" = 45958
    "D54" = 180.28
    "E54" = 174.75
    "F54" = 184.75
    "A55" = 45958
    "D55" = 167.93
    "E55" = 171.98
    "F55" = 181.98
    "A56" = 45958
    "D56" = 170.2
    "A57" = 45958
    "D57" = 169.88
    "E57" = 166.25
    "A58" = 45958
    "D58" = 165.78
    "E58" = 162.30000000000001
    "F58" = 172.3
    "A59" = 45958
    "D59" = 172.61
    "E59" = 172.96
    "A60" = 45955
    "D60" = 177.81
    "E60" = 173.44
    "F60" = 183.44
    "A61" = 45955
    "D61" = 165.46
    "E61" = 170.77
    "F61" = 180.77
    "A62" = 45955
    "D62" = 167.74
    "A63" = 45955
    "D63" = 167.42
    "E63" = 165.04
    "A64" = 45955
    "D64" = 163.32
    "E64" = 161.09
    "F64" = 171.09
    "A65" = 45955
    "D65" = 170.15
    "E65" = 171.65
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}